$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# Rename header cells: Flag -> Run, sno -> Sno
$ws.Range("B1").Value = "Run"
$ws.Range("A1").Value = "Sno"

# Update the Flag/Run column values
$ws.Range("B2").Value = "no"
$ws.Range("B3").Value = "yes"
$ws.Range("B4").Value = ""
$ws.Range("B5").Value = "yes"

# Hide Sheet2 and Sheet3
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Visible = $false
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Visible = $false
